$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update the "Förändrad" date column (C) for rows 2-27 from 45305 to 45306
for ($row = 2; $row -le 27; $row++) {
    $ws.Cells.Item($row, 3).Value = 45306
}
